$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of results (rows 69-78), columns C, F, G only (others left blank)
$data = @(
    @{ Row = 69; C = 388016;  F = 24.029699999999998; G = 4 },
    @{ Row = 70; C = 225017;  F = 31.9802;             G = 3 },
    @{ Row = 71; C = 113044;  F = 19.010899999999999; G = 3 },
    @{ Row = 72; C = "12375764_10154354426419428_2121622626_o"; F = 23.587499999999999; G = 4 },
    @{ Row = 73; C = "fabric";      F = 20.681799999999999; G = 6 },
    @{ Row = 74; C = "10472953_992127864183797_1990666493_n";  F = 9.7666000000000004; G = 4 },
    @{ Row = 75; C = "linfocitos1"; F = 5.0885999999999996;  G = 12 },
    @{ Row = 76; C = 135069;  F = 12.0886;               G = 7 },
    @{ Row = 77; C = "leucolinf";   F = 4.7206999999999999;  G = 18 },
    @{ Row = 78; C = "showimage";   F = 21.145600000000002; G = 4 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

# Update the view to reflect where the user scrolled/selected after adding rows
$ws.Application.ActiveWindow.ScrollRow = 73
$ws.Range("E84").Select()
